$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2168.19025045373
$ws.Range("C2").Value = 2167.314610709427
$ws.Range("D2").Value = 2165.691048665342
$ws.Range("E2").Value = 2160.906096720676
$ws.Range("B3").Value = 2183.438124058833
$ws.Range("C3").Value = 2183.632415180711
$ws.Range("D3").Value = 2183.136329287854
$ws.Range("E3").Value = 2178.466453297492
$ws.Range("B4").Value = 2134.223707206921
$ws.Range("C4").Value = 2131.347061543535
$ws.Range("D4").Value = 2127.598171040406
$ws.Range("E4").Value = 2122.310463277523
$ws.Range("B5").Value = 2184.588446198372
$ws.Range("C5").Value = 2185.449126116547
$ws.Range("D5").Value = 2185.201902321746
$ws.Range("E5").Value = 2180.197948453717
$ws.Range("B6").Value = 2190.810902805146
$ws.Range("C6").Value = 2191.656695839558
$ws.Range("D6").Value = 2191.64712965028
$ws.Range("E6").Value = 2186.955418566179
$ws.Range("B7").Value = 2172.524713790897
$ws.Range("C7").Value = 2172.86595052552
$ws.Range("D7").Value = 2171.940903516458
$ws.Range("E7").Value = 2166.472776123554
$ws.Range("B8").Value = 2174.232017974375
$ws.Range("C8").Value = 2173.684592583933
$ws.Range("D8").Value = 2172.503856724016
$ws.Range("E8").Value = 2167.740513980736
$ws.Range("B9").Value = 2186.103315481719
$ws.Range("C9").Value = 2187.12321845688
$ws.Range("D9").Value = 2187.02243935776
$ws.Range("E9").Value = 2181.940540741372
$ws.Range("B10").Value = 1946.701420854897
$ws.Range("C10").Value = 1937.729501449272
$ws.Range("D10").Value = 1922.900717917014
$ws.Range("E10").Value = 1908.84393146966
$ws.Range("B11").Value = 1895.566769401754
$ws.Range("C11").Value = 1884.386919791947
$ws.Range("D11").Value = 1866.433802825766
$ws.Range("E11").Value = 1850.571204774973
$ws.Range("B12").Value = 1637.698552672618
$ws.Range("C12").Value = 1614.595651276729
$ws.Range("D12").Value = 1582.176651645972
$ws.Range("E12").Value = 1557.822975848908
$ws.Range("B13").Value = 1929.418281600073
$ws.Range("C13").Value = 1919.499697099683
$ws.Range("D13").Value = 1903.651813317016
$ws.Range("E13").Value = 1889.092437598045
